# paises.xlsx data refresh: "Update countries & provincias Spain"
# Updates daily COVID-19 stats, inserts "Burkina Faso" above "Reunion",
# moves "Montenegro" above "Ghana", and bumps the "updated at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------
# Simple per-row statistic refreshes (no rows inserted/removed)
# ---------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 323953
$ws.Range("C4").Value = 12596
$ws.Range("D4").Value = 16598
$ws.Range("E4").Value = 298170
$ws.Range("G4").Value = 733
$ws.Range("H4").Value = 9185

# Row 7 - Alemania
$ws.Range("B7").Value = 98765
$ws.Range("C7").Value = 2673
$ws.Range("E7").Value = 70841
$ws.Range("G7").Value = 80
$ws.Range("H7").Value = 1524

# Row 41
$ws.Range("D41").Value = 793
$ws.Range("E41").Value = 1353

# Row 52
$ws.Range("D52").Value = 280
$ws.Range("E52").Value = 1127

# Row 56
$ws.Range("F56").Value = 25

# Row 74
$ws.Range("B74").Value = 584
$ws.Range("C74").Value = 53
$ws.Range("D74").Value = 42
$ws.Range("E74").Value = 536

# Row 110 - Georgia
$ws.Range("B110").Value = 174
$ws.Range("C110").Value = 12
$ws.Range("E110").Value = 136

# ---------------------------------------------------------------
# "Burkina Faso" climbs above "Reunion" (row 90), old Burkina Faso
# row (which shifts to row 95 once the new row is inserted) is
# then removed since it is now a duplicate entry.
# ---------------------------------------------------------------

$ws.Rows(90).Insert()
$ws.Range("A90").Value = "Burkina Faso"
$ws.Range("B90").Value = 345
$ws.Range("C90").Value = 27
$ws.Range("D90").Value = 90
$ws.Range("E90").Value = 238
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 17
$ws.Rows(95).Delete()

# ---------------------------------------------------------------
# "Montenegro" climbs above "Ghana" (row 106), old Montenegro row
# (which shifts to row 108 once the new row is inserted) is then
# removed since it is now a duplicate entry.
# ---------------------------------------------------------------

$ws.Rows(106).Insert()
$ws.Range("A106").Value = "Montenegro"
$ws.Range("B106").Value = 214
$ws.Range("C106").Value = 13
$ws.Range("D106").Value = 1
$ws.Range("E106").Value = 211
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2
$ws.Rows(108).Delete()

# ---------------------------------------------------------------
# Header timestamp ("Datos actualizados a ...")
# ---------------------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 18:52"
